$wb = $excel.ActiveWorkbook

# 1. Set Sistemas!F19 to "marce" (new shared string inserted first, per commit ordering)
$sistemas = $wb.Worksheets.Item("Sistemas")
$sistemas.Range("F19").Value = "marce"

# 2. Update selection on Sistemas (bottomRight pane) from C22 to H16
$sistemas.Range("H16").Select()

# 3. Update selection on Usuarios from D15 to C15
$usuarios = $wb.Worksheets.Item("Usuarios")
$usuarios.Range("C15").Select()

# 4. Add the new "Tarjetas" sheet at the end of the workbook
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$tarjetas = $wb.Worksheets.Add($null, $lastSheet)
$tarjetas.Name = "Tarjetas"

# Column widths (closest achievable given host pixel quantization)
$tarjetas.Columns.Item(1).ColumnWidth = 29.833333333333332
$tarjetas.Columns.Item(3).ColumnWidth = 20.333333333333332
$tarjetas.Columns.Item(4).ColumnWidth = 23.833333333333332

# Header row
$tarjetas.Cells.Item(1,1).Value = "Nombre y apellido"
$tarjetas.Cells.Item(1,2).Value = "DNI"
$tarjetas.Cells.Item(1,3).Value = "Área"
$tarjetas.Cells.Item(1,4).Value = "Número de tarjeta"

# Column A (names) filled top-to-bottom first, matching shared-string insertion order
$tarjetas.Cells.Item(2,1).Value = "Gabriel Sánchez"
$tarjetas.Cells.Item(3,1).Value = "Ricardo Hurtado"
$tarjetas.Cells.Item(4,1).Value = "Germán Kutalek"
$tarjetas.Cells.Item(5,1).Value = "Santiago Franceschi"
$tarjetas.Cells.Item(6,1).Value = "Oscar Faranna"
$tarjetas.Cells.Item(7,1).Value = "Wilson Gómez"
$tarjetas.Cells.Item(8,1).Value = "Jorge Balcazar"
$tarjetas.Cells.Item(9,1).Value = "Enzo Castillo"
$tarjetas.Cells.Item(10,1).Value = "Manuel Umpierrez"
$tarjetas.Cells.Item(11,1).Value = "Jonatan Soto"
$tarjetas.Cells.Item(12,1).Value = "Adrián Meana"
$tarjetas.Cells.Item(13,1).Value = "Gabriel Zeñiuk"
$tarjetas.Cells.Item(14,1).Value = "Gustavo Luis Carca"
$tarjetas.Cells.Item(15,1).Value = "Sergio Ochoa"
$tarjetas.Cells.Item(16,1).Value = "Sergio Eduardo Gutierre"
$tarjetas.Cells.Item(17,1).Value = "Miguel Ángel Riccitelli"
$tarjetas.Cells.Item(18,1).Value = "Mirta Rojas"
$tarjetas.Cells.Item(19,1).Value = "Llavero 01"
$tarjetas.Cells.Item(20,1).Value = "Llavero 02"
$tarjetas.Cells.Item(21,1).Value = "Randy Speake"
$tarjetas.Cells.Item(22,1).Value = "Pablo Maccari"
$tarjetas.Cells.Item(23,1).Value = "Mantenimiento (master)"
$tarjetas.Cells.Item(24,1).Value = "Bombero 01 (Fernando Sánchez)"

# Column B (DNI numbers)
$tarjetas.Cells.Item(2,2).Value = 23426675
$tarjetas.Cells.Item(3,2).Value = 35897670
$tarjetas.Cells.Item(4,2).Value = 37878033
$tarjetas.Cells.Item(5,2).Value = 14972197
$tarjetas.Cells.Item(6,2).Value = 28746454
$tarjetas.Cells.Item(7,2).Value = 94929982
$tarjetas.Cells.Item(8,2).Value = 37900747
$tarjetas.Cells.Item(9,2).Value = 34932264
$tarjetas.Cells.Item(10,2).Value = 92821036
$tarjetas.Cells.Item(11,2).Value = 32604554
$tarjetas.Cells.Item(12,2).Value = 29704359
$tarjetas.Cells.Item(13,2).Value = 27067242
$tarjetas.Cells.Item(14,2).Value = 11293948
$tarjetas.Cells.Item(15,2).Value = 16965144
$tarjetas.Cells.Item(16,2).Value = 24222400
$tarjetas.Cells.Item(17,2).Value = 27257335
$tarjetas.Cells.Item(18,2).Value = 20640391
$tarjetas.Cells.Item(22,2).Value = 23179837
$tarjetas.Cells.Item(24,2).Value = 32554035

# Column C (Área) filled top-to-bottom
$tarjetas.Cells.Item(2,3).Value = "Mantenimiento"
$tarjetas.Cells.Item(3,3).Value = "Mantenimiento"
$tarjetas.Cells.Item(4,3).Value = "Mantenimiento"
$tarjetas.Cells.Item(5,3).Value = "Mantenimiento"
$tarjetas.Cells.Item(6,3).Value = "Mantenimiento"
$tarjetas.Cells.Item(7,3).Value = "Mantenimiento"
$tarjetas.Cells.Item(8,3).Value = "Mantenimiento"
$tarjetas.Cells.Item(9,3).Value = "Mantenimiento"
$tarjetas.Cells.Item(10,3).Value = "CTO"
$tarjetas.Cells.Item(11,3).Value = "Operativa y Logística"
$tarjetas.Cells.Item(12,3).Value = "Operativa y Logística"
$tarjetas.Cells.Item(13,3).Value = "Limpieza"
$tarjetas.Cells.Item(14,3).Value = "Seguridad"
$tarjetas.Cells.Item(15,3).Value = "Seguridad"
$tarjetas.Cells.Item(16,3).Value = "Seguridad"
$tarjetas.Cells.Item(17,3).Value = "Seguridad"
$tarjetas.Cells.Item(18,3).Value = "Seguridad"
$tarjetas.Cells.Item(19,3).Value = "Espacios Físicos"
$tarjetas.Cells.Item(20,3).Value = "Espacios Físicos"
$tarjetas.Cells.Item(21,3).Value = "Espacios Físicos"
$tarjetas.Cells.Item(22,3).Value = "Control y Planificación"
$tarjetas.Cells.Item(24,3).Value = "Control y Planificación"

# Column D (Número de tarjeta) filled top-to-bottom
$tarjetas.Cells.Item(2,4).Value = "0008382400 | 127,59328"
$tarjetas.Cells.Item(3,4).Value = "0008379433 | 127,56361"
$tarjetas.Cells.Item(4,4).Value = "0008382623 | 127,59551"
$tarjetas.Cells.Item(5,4).Value = "0008383530 | 127,60458"
$tarjetas.Cells.Item(6,4).Value = "0008390124 | 128,01516"
$tarjetas.Cells.Item(7,4).Value = "0008388805 | 128,00197"
$tarjetas.Cells.Item(8,4).Value = "0008384610 | 127,61538"
$tarjetas.Cells.Item(9,4).Value = "0008385792 | 127,62720"
$tarjetas.Cells.Item(10,4).Value = "0008389945 | 128,01337"
$tarjetas.Cells.Item(11,4).Value = "0008389664 | 128,01056"
$tarjetas.Cells.Item(12,4).Value = "0008391207 | 128,02599"
$tarjetas.Cells.Item(13,4).Value = "0008388946 | 128,00338"
$tarjetas.Cells.Item(14,4).Value = "0008386875 | 127,63803"
$tarjetas.Cells.Item(15,4).Value = "0008387294 | 127,64222"
$tarjetas.Cells.Item(16,4).Value = "0008380931 | 127,57859"
$tarjetas.Cells.Item(17,4).Value = "0008383258 | 127,60186"
$tarjetas.Cells.Item(18,4).Value = "0008380357 | 127,57285"
$tarjetas.Cells.Item(19,4).Value = "0008386494 | 127,63422"
$tarjetas.Cells.Item(20,4).Value = "0008385642 | 127,62570"
$tarjetas.Cells.Item(21,4).Value = "0008381101 | 127,58029"
$tarjetas.Cells.Item(22,4).Value = "0008382363 | 127,59291"
$tarjetas.Cells.Item(23,4).Value = "0008380237 | 127,57165"
$tarjetas.Cells.Item(24,4).Value = "0008382231 | 127,59159"

# Final selection/active cell on the new sheet
$tarjetas.Range("D25").Select()

Write-Output "done"
